$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 12:44"

# Row 14 (Iran) - updated case numbers
$ws.Range("B14").Value = 271606
$ws.Range("C14").Value = 2166
$ws.Range("D14").Value = 235300
$ws.Range("E14").Value = 22327
$ws.Range("G14").Value = 188
$ws.Range("H14").Value = 13979

# Row 25 (Catar) - updated case numbers
$ws.Range("B25").Value = 106308
$ws.Range("C25").Value = 410
$ws.Range("D25").Value = 103023
$ws.Range("E25").Value = 3131
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 154

# Row 39 (Emiratos Arabes Unidos) - updated case numbers
$ws.Range("B39").Value = 56711
$ws.Range("C39").Value = 289
$ws.Range("D39").Value = 48917
$ws.Range("E39").Value = 7456
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 338

# Row 48 (Rumania) - updated case numbers
$ws.Range("B48").Value = 36691
$ws.Range("C48").Value = 889
$ws.Range("D48").Value = 22488
$ws.Range("E48").Value = 12194
$ws.Range("G48").Value = 21
$ws.Range("H48").Value = 2009

# Rows 74/75: El Salvador overtakes Australia in the ranking, so the two
# countries swap places. Row 74 becomes El Salvador (with new, larger
# numbers) and row 75 becomes Australia (keeping its previous numbers).
$ws.Range("A74").Value = "El Salvador"
$ws.Range("B74").Value = 11508
$ws.Range("C74").Value = 301
$ws.Range("D74").Value = 6560
$ws.Range("E74").Value = 4624
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 15
$ws.Range("H74").Value = 324

$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 11438
$ws.Range("C75").Value = 203
$ws.Range("D75").Value = 8158
$ws.Range("E75").Value = 3162
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 118

# Row 126 (Hong Kong) - updated case numbers
$ws.Range("B126").Value = 1778
$ws.Range("C126").Value = 64
$ws.Range("D126").Value = 1274
$ws.Range("E126").Value = 492
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 12
